$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on percentage-like cells in column H so values such as "55%" are
# stored as literal text (matching the source data) instead of being auto-converted
# to a numeric percentage by Excel.
$ws.Range('H4').NumberFormat = "@"
$ws.Range('H12').NumberFormat = "@"
$ws.Range('H13').NumberFormat = "@"
$ws.Range('H14').NumberFormat = "@"
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H22').NumberFormat = "@"
$ws.Range('H23').NumberFormat = "@"
$ws.Range('H25').NumberFormat = "@"
$ws.Range('H26').NumberFormat = "@"
$ws.Range('H28').NumberFormat = "@"
$ws.Range('H30').NumberFormat = "@"
$ws.Range('H34').NumberFormat = "@"
$ws.Range('H37').NumberFormat = "@"
$ws.Range('H38').NumberFormat = "@"
$ws.Range('H39').NumberFormat = "@"
$ws.Range('H40').NumberFormat = "@"
$ws.Range('H46').NumberFormat = "@"

$ws.Range('E2').Value = '2026-02-20 06:48:10'
$ws.Range('E3').Value = '2026-02-20 06:48:13'
$ws.Range('N3').Value = '-6.3 °C 6:24 TU'
$ws.Range('E4').Value = '2026-02-20 06:48:15'
$ws.Range('H4').Value = '55%'
$ws.Range('J4').Value = '1018.9 hPa'
$ws.Range('N4').Value = '5.2 °C 6:03 TU'
$ws.Range('O4').Value = '8.3 °C'
$ws.Range('E5').Value = '2026-02-20 06:48:17'
$ws.Range('I5').Value = '1.2 mm'
$ws.Range('E6').Value = '2026-02-20 06:48:20'
$ws.Range('J6').Value = '1019.0 hPa'
$ws.Range('N6').Value = '3.0 °C 6:09 TU'
$ws.Range('O6').Value = '5.2 °C'
$ws.Range('E7').Value = '2026-02-20 06:48:22'
$ws.Range('J7').Value = '1018.5 hPa'
$ws.Range('E8').Value = '2026-02-20 06:48:24'
$ws.Range('J8').Value = '1019.4 hPa'
$ws.Range('M8').Value = '7.9 °C 6:24 TU'
$ws.Range('O8').Value = '7.0 °C'
$ws.Range('E9').Value = '2026-02-20 06:48:27'
$ws.Range('O9').Value = '12.4 °C'
$ws.Range('E10').Value = '2026-02-20 06:48:29'
$ws.Range('E11').Value = '2026-02-20 06:48:32'
$ws.Range('E12').Value = '2026-02-20 06:48:34'
$ws.Range('H12').Value = '51%'
$ws.Range('O12').Value = '12.1 °C'
$ws.Range('E13').Value = '2026-02-20 06:48:36'
$ws.Range('H13').Value = '47%'
$ws.Range('J13').Value = '1020.4 hPa'
$ws.Range('N13').Value = '1.9 °C 6:11 TU'
$ws.Range('O13').Value = '4.8 °C'
$ws.Range('E14').Value = '2026-02-20 06:48:39'
$ws.Range('H14').Value = '54%'
$ws.Range('N14').Value = '9.3 °C 6:14 TU'
$ws.Range('O14').Value = '10.1 °C'
$ws.Range('E15').Value = '2026-02-20 06:48:41'
$ws.Range('N15').Value = '11.5 °C 6:29 TU'
$ws.Range('O15').Value = '12.6 °C'
$ws.Range('E16').Value = '2026-02-20 06:48:43'
$ws.Range('H16').Value = '56%'
$ws.Range('E17').Value = '2026-02-20 06:48:46'
$ws.Range('E18').Value = '2026-02-20 06:48:48'
$ws.Range('J18').Value = '1019.3 hPa'
$ws.Range('N18').Value = '-0.1 °C 6:17 TU'
$ws.Range('O18').Value = '1.6 °C'
$ws.Range('E19').Value = '2026-02-20 06:48:51'
$ws.Range('E20').Value = '2026-02-20 06:48:53'
$ws.Range('N20').Value = '-6.2 °C 6:04 TU'
$ws.Range('E21').Value = '2026-02-20 06:48:55'
$ws.Range('J21').Value = '1020.5 hPa'
$ws.Range('N21').Value = '3.4 °C 6:13 TU'
$ws.Range('O21').Value = '6.1 °C'
$ws.Range('E22').Value = '2026-02-20 06:48:58'
$ws.Range('H22').Value = '61%'
$ws.Range('M22').Value = '-4.2 °C 6:15 TU'
$ws.Range('O22').Value = '-6.4 °C'
$ws.Range('E23').Value = '2026-02-20 06:49:00'
$ws.Range('H23').Value = '84%'
$ws.Range('I23').Value = '3.6 mm'
$ws.Range('N23').Value = '-7.1 °C 6:01 TU'
$ws.Range('E24').Value = '2026-02-20 06:49:02'
$ws.Range('J24').Value = '1023.2 hPa'
$ws.Range('N24').Value = '4.9 °C 6:28 TU'
$ws.Range('O24').Value = '6.1 °C'
$ws.Range('E25').Value = '2026-02-20 06:49:05'
$ws.Range('H25').Value = '70%'
$ws.Range('I25').Value = '4.5 mm'
$ws.Range('M25').Value = '-3.4 °C 6:26 TU'
$ws.Range('E26').Value = '2026-02-20 06:49:07'
$ws.Range('H26').Value = '41%'
$ws.Range('J26').Value = '1018.7 hPa'
$ws.Range('N26').Value = '2.0 °C 6:03 TU'
$ws.Range('E27').Value = '2026-02-20 06:49:10'
$ws.Range('E28').Value = '2026-02-20 06:49:12'
$ws.Range('H28').Value = '75%'
$ws.Range('J28').Value = '1020.0 hPa'
$ws.Range('N28').Value = '0.2 °C 6:29 TU'
$ws.Range('O28').Value = '2.4 °C'
$ws.Range('E29').Value = '2026-02-20 06:49:15'
$ws.Range('E30').Value = '2026-02-20 06:49:17'
$ws.Range('H30').Value = '66%'
$ws.Range('J30').Value = '1018.5 hPa'
$ws.Range('O30').Value = '8.6 °C'
$ws.Range('E31').Value = '2026-02-20 06:49:19'
$ws.Range('J31').Value = '1017.0 hPa'
$ws.Range('N31').Value = '9.3 °C 6:25 TU'
$ws.Range('O31').Value = '10.4 °C'
$ws.Range('E32').Value = '2026-02-20 06:49:22'
$ws.Range('N32').Value = '1.2 °C 6:29 TU'
$ws.Range('O32').Value = '1.8 °C'
$ws.Range('E33').Value = '2026-02-20 06:49:24'
$ws.Range('J33').Value = '1019.5 hPa'
$ws.Range('O33').Value = '4.4 °C'
$ws.Range('E34').Value = '2026-02-20 06:49:27'
$ws.Range('H34').Value = '62%'
$ws.Range('M34').Value = '-0.8 °C 6:18 TU'
$ws.Range('O34').Value = '-2.0 °C'
$ws.Range('E35').Value = '2026-02-20 06:49:29'
$ws.Range('J35').Value = '1024.2 hPa'
$ws.Range('E36').Value = '2026-02-20 06:49:31'
$ws.Range('J36').Value = '1018.8 hPa'
$ws.Range('E37').Value = '2026-02-20 06:49:34'
$ws.Range('H37').Value = '68%'
$ws.Range('J37').Value = '1021.6 hPa'
$ws.Range('N37').Value = '-0.8 °C 6:29 TU'
$ws.Range('O37').Value = '1.7 °C'
$ws.Range('E38').Value = '2026-02-20 06:49:36'
$ws.Range('H38').Value = '76%'
$ws.Range('E39').Value = '2026-02-20 06:49:39'
$ws.Range('G39').Value = '55 cm'
$ws.Range('H39').Value = '73%'
$ws.Range('M39').Value = '-4.3 °C 6:23 TU'
$ws.Range('O39').Value = '-5.8 °C'
$ws.Range('E40').Value = '2026-02-20 06:49:41'
$ws.Range('H40').Value = '47%'
$ws.Range('J40').Value = '1021.1 hPa'
$ws.Range('O40').Value = '7.5 °C'
$ws.Range('E41').Value = '2026-02-20 06:49:43'
$ws.Range('J41').Value = '1020.5 hPa'
$ws.Range('N41').Value = '10.0 °C 6:13 TU'
$ws.Range('E42').Value = '2026-02-20 06:49:46'
$ws.Range('O42').Value = '4.4 °C'
$ws.Range('E43').Value = '2026-02-20 06:49:48'
$ws.Range('E44').Value = '2026-02-20 06:49:50'
$ws.Range('I44').Value = '4.5 mm'
$ws.Range('E45').Value = '2026-02-20 06:49:53'
$ws.Range('J45').Value = '1027.8 hPa'
$ws.Range('N45').Value = '1.6 °C 6:15 TU'
$ws.Range('E46').Value = '2026-02-20 06:49:55'
$ws.Range('H46').Value = '60%'
$ws.Range('J46').Value = '1024.1 hPa'
